$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Persona" row (row 2). Deleting the row shifts
# "Organización" and "Denuncia" up, matching the target layout:
#   Row 2: Organización
#   Row 3: Denuncia
$ws.Rows.Item(2).Delete()
